$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 177.4
$ws.Range("I6").Value = 177.4
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 532.2
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -420.2
$ws.Range("N6").ClearContents()
# Row 51
$ws.Range("H51").Value = 47249.5
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 47249.5
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 47249.5
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -48217.5
# Row 87
$ws.Range("H87").Value = 98281.75
$ws.Range("I87").Value = 20000
$ws.Range("J87").Value = 124375.664
$ws.Range("K87").Value = 20000
$ws.Range("L87").Value = 124375.664
$ws.Range("M87").Value = -18752
$ws.Range("N87").Value = -126871.664
# Row 90
$ws.Range("H90").Value = 98281.75
$ws.Range("I90").Value = 20000
$ws.Range("J90").Value = 124375.664
$ws.Range("K90").Value = 60000
$ws.Range("L90").Value = 373126.992
$ws.Range("M90").Value = -53760
$ws.Range("N90").Value = -385606.992
# Row 96
$ws.Range("H96").Value = 1870.4
$ws.Range("J96").Value = 2724.6667
$ws.Range("L96").Value = 8174.000100000001
$ws.Range("N96").Value = -10920.0001
# Row 111
$ws.Range("H111").Value = 11339.343
$ws.Range("I111").Value = 9723.294
$ws.Range("J111").Value = 12865.611
$ws.Range("K111").Value = 29169.882
$ws.Range("L111").Value = 38596.833
$ws.Range("M111").Value = -26102.882
$ws.Range("N111").Value = -44730.833
# Row 132
$ws.Range("H132").Value = 2788.804
$ws.Range("I132").Value = 2796.2856
$ws.Range("J132").Value = 2605.5
$ws.Range("K132").Value = 8388.856800000001
$ws.Range("L132").Value = 7816.5
$ws.Range("M132").Value = -5858.856800000001
$ws.Range("N132").Value = -12876.5
# Row 137
$ws.Range("H137").Value = 1868.28
$ws.Range("I137").Value = 1717.75
$ws.Range("J137").Value = 2135.889
$ws.Range("K137").Value = 5153.25
$ws.Range("L137").Value = 6407.667
$ws.Range("M137").Value = -2603.25
$ws.Range("N137").Value = -11507.667

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 3190.6365
$ws.Range("I2").Value = 2085.4285
$ws.Range("J2").Value = 5124.75
$ws.Range("K2").Value = 2085.4285
$ws.Range("L2").Value = 5124.75
$ws.Range("M2").Value = -1972.4285
$ws.Range("N2").Value = -5350.75
# Row 55
$ws.Range("H55").Value = 14999
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
# Row 61
$ws.Range("H61").Value = 4748.4
$ws.Range("J61").Value = 3000
$ws.Range("L61").Value = 3000
$ws.Range("N61").Value = -3424
# Row 116
$ws.Range("H116").Value = 3190.6365
$ws.Range("I116").Value = 2085.4285
$ws.Range("J116").Value = 5124.75
$ws.Range("K116").Value = 2085.4285
$ws.Range("L116").Value = 5124.75
$ws.Range("M116").Value = 208.5715
$ws.Range("N116").Value = -9712.75
# Row 136
$ws.Range("H136").Value = 4748.4
$ws.Range("J136").Value = 3000
$ws.Range("L136").Value = 9000
$ws.Range("N136").Value = -14100

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 3190.6365
$ws.Range("I3").Value = 2085.4285
$ws.Range("J3").Value = 5124.75
$ws.Range("K3").Value = 2085.4285
$ws.Range("L3").Value = 5124.75
$ws.Range("M3").Value = -1971.4285
$ws.Range("N3").Value = -5352.75

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 709.13336
$ws.Range("I22").Value = 423.8
$ws.Range("J22").Value = 1279.8
$ws.Range("K22").Value = 423.8
$ws.Range("L22").Value = 1279.8
$ws.Range("M22").Value = -73.80000000000001
$ws.Range("N22").Value = -1979.8
# Row 31
$ws.Range("H31").Value = 4103.037
$ws.Range("I31").Value = 4198.4
$ws.Range("J31").Value = 3830.5715
$ws.Range("K31").Value = 4198.4
$ws.Range("L31").Value = 3830.5715
$ws.Range("M31").Value = -3903.4
$ws.Range("N31").Value = -4420.5715
# Row 34
$ws.Range("H34").Value = 4103.037
$ws.Range("I34").Value = 4198.4
$ws.Range("J34").Value = 3830.5715
$ws.Range("K34").Value = 4198.4
$ws.Range("L34").Value = 3830.5715
$ws.Range("M34").Value = -3996.4
$ws.Range("N34").Value = -4234.5715
# Row 58
$ws.Range("H58").Value = 5287.3335
$ws.Range("I58").Value = 5581.294
$ws.Range("J58").Value = 4038
$ws.Range("K58").Value = 5581.294
$ws.Range("L58").Value = 4038
$ws.Range("M58").Value = -5378.294
$ws.Range("N58").Value = -4444
# Row 125
$ws.Range("H125").Value = 89998
$ws.Range("J125").Value = 89998
$ws.Range("L125").Value = 89998
$ws.Range("N125").Value = -94918
# Row 134
$ws.Range("H134").Value = 7783.5537
$ws.Range("I134").Value = 6910.4346
$ws.Range("J134").Value = 11799.9
$ws.Range("K134").Value = 20731.3038
$ws.Range("L134").Value = 35399.7
$ws.Range("M134").Value = -18196.3038
$ws.Range("N134").Value = -40469.7
# Row 136
$ws.Range("H136").Value = 5287.3335
$ws.Range("I136").Value = 5581.294
$ws.Range("J136").Value = 4038
$ws.Range("K136").Value = 16743.882
$ws.Range("L136").Value = 12114
$ws.Range("M136").Value = -14193.882
$ws.Range("N136").Value = -17214

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 114
$ws.Range("H114").Value = 1826.875
$ws.Range("J114").Value = 2131.1667
$ws.Range("L114").Value = 6393.500100000001
$ws.Range("N114").Value = -12901.5001
# Row 137
$ws.Range("H137").Value = 3276.6667
$ws.Range("I137").Value = 3784.4285
$ws.Range("K137").Value = 11353.2855
$ws.Range("M137").Value = -6253.2855

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 3999
$ws.Range("I22").Value = 3999
$ws.Range("K22").Value = 3999
$ws.Range("M22").Value = -3704
# Row 27
$ws.Range("H27").Value = 3999
$ws.Range("I27").Value = 3999
$ws.Range("K27").Value = 3999
$ws.Range("M27").Value = -3892
# Row 46
$ws.Range("H46").Value = 5419.615
$ws.Range("I46").Value = 6183
$ws.Range("K46").Value = 6183
$ws.Range("M46").Value = -5995
# Row 55
$ws.Range("H55").Value = 281.66666
$ws.Range("I55").Value = 608.75
$ws.Range("K55").Value = 608.75
$ws.Range("M55").Value = -435.75
# Row 114
$ws.Range("H114").Value = 88997.336
$ws.Range("J114").Value = 88997.336
$ws.Range("L114").Value = 88997.336
$ws.Range("N114").Value = -97675.336
# Row 120
$ws.Range("H120").Value = 73329
$ws.Range("J120").Value = 73329
$ws.Range("L120").Value = 73329
$ws.Range("N120").Value = -83005
# Row 136
$ws.Range("H136").Value = 1897.2632
$ws.Range("I136").Value = 1871.6666
$ws.Range("K136").Value = 5614.9998
$ws.Range("M136").Value = -3064.9998

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 21
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
# Row 28
$ws.Range("H28").Value = 20000
$ws.Range("I28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("M28").ClearContents()
# Row 30
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()
# Row 35
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
# Row 41
$ws.Range("H41").Value = 12000
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 12000
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 12000
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -12780
# Row 101
$ws.Range("H101").Value = 24683.572
$ws.Range("J101").Value = 18464.166
$ws.Range("L101").Value = 18464.166
$ws.Range("N101").Value = -24954.166
# Row 136
$ws.Range("H136").Value = 3949.2285
$ws.Range("I136").Value = 3869.8965
$ws.Range("K136").Value = 11609.6895
$ws.Range("M136").Value = -9059.6895
